$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report B440")

# Carry-over balance feeding the running Progress column (G6 = G5 + E6, etc.)
$ws.Cells.Item(5, 7).Value = -2419

# Plan (D) and Actual (E) values per row, as described by the diff.
# Rows without an explicit key keep their existing values (weekend rows with
# no Plan/Actual entries, row 36 blank separator, etc.)
$data = @{
    6  = @{ D = 110; E = 117 }
    7  = @{ D = 110; E = 114 }
    8  = @{ D = 110; E = 106 }
    9  = @{ D = 110; E = 119 }
    12 = @{ D = 110; E = 125 }
    13 = @{ D = 144; E = 120 }
    14 = @{ D = 110; E = 110 }
    15 = @{ D = 110; E = 104 }
    16 = @{ D = 110; E = 114 }
    19 = @{ D = 110; E = 122 }
    20 = @{ D = 144; E = 148 }
    21 = @{ D = 110; E = 114 }
    22 = @{ D = 110; E = 135 }
    23 = @{ D = 110; E = 87 }
    26 = @{ D = 110 }
    27 = @{ D = 110 }
    28 = @{ D = 110 }
    29 = @{ D = 110 }
    30 = @{ D = 110 }
    33 = @{ D = 110 }
    34 = @{ D = 110 }
    35 = @{ D = 41 }
}

foreach ($row in $data.Keys) {
    $entry = $data[$row]
    if ($entry.ContainsKey("D")) {
        $ws.Cells.Item($row, 4).Value = $entry.D
    }
    if ($entry.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $entry.E
    }
}

$wb.Application.CalculateFull()
